# Fix significance analysis bug: num_X columns incorrectly defaulted to 5 in many
# cases. Re-output corrected num_X counts and recompute p_X = num_X / p_cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N3").Value = 0.2727272727272727
$ws.Range("O3").Value = 3.0
$ws.Range("L4").Value = 0.08333333333333333
$ws.Range("M4").Value = 1.0
$ws.Range("N4").Value = 0.08333333333333333
$ws.Range("O4").Value = 1.0
$ws.Range("D5").Value = 0.07142857142857142
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 0.09523809523809523
$ws.Range("G5").Value = 4.0
$ws.Range("J5").Value = 0.023809523809523808
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.07142857142857142
$ws.Range("M5").Value = 3.0
$ws.Range("H6").Value = 0.037037037037037035
$ws.Range("I6").Value = 1.0
$ws.Range("J6").Value = 0.037037037037037035
$ws.Range("K6").Value = 1.0
$ws.Range("L6").Value = 0.07407407407407407
$ws.Range("M6").Value = 2.0
$ws.Range("N6").Value = 0.07407407407407407
$ws.Range("O6").Value = 2.0
$ws.Range("D7").Value = 0.07692307692307693
$ws.Range("E7").Value = 1.0
$ws.Range("F7").Value = 0.15384615384615385
$ws.Range("G7").Value = 2.0
$ws.Range("H7").Value = 0.15384615384615385
$ws.Range("I7").Value = 2.0
$ws.Range("L7").Value = 0.07692307692307693
$ws.Range("M7").Value = 1.0
$ws.Range("N7").Value = 0.15384615384615385
$ws.Range("O7").Value = 2.0
$ws.Range("D8").Value = 0.018867924528301886
$ws.Range("E8").Value = 1.0
$ws.Range("F8").Value = 0.07547169811320754
$ws.Range("G8").Value = 4.0
$ws.Range("J8").Value = 0.018867924528301886
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.07547169811320754
$ws.Range("M8").Value = 4.0
$ws.Range("D9").Value = 0.03333333333333333
$ws.Range("E9").Value = 1.0
$ws.Range("F9").Value = 0.03333333333333333
$ws.Range("G9").Value = 1.0
$ws.Range("N9").Value = 0.1
$ws.Range("O9").Value = 3.0
$ws.Range("F10").Value = 0.05128205128205128
$ws.Range("G10").Value = 2.0
$ws.Range("H10").Value = 0.10256410256410256
$ws.Range("I10").Value = 4.0
$ws.Range("J10").Value = 0.05128205128205128
$ws.Range("K10").Value = 2.0
$ws.Range("L10").Value = 0.05128205128205128
$ws.Range("M10").Value = 2.0
$ws.Range("N10").Value = 0.07692307692307693
$ws.Range("O10").Value = 3.0
$ws.Range("D11").Value = 0.2
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 0.2
$ws.Range("G11").Value = 3.0
$ws.Range("J11").Value = 0.06666666666666667
$ws.Range("K11").Value = 1.0
$ws.Range("L11").Value = 0.06666666666666667
$ws.Range("M11").Value = 1.0
$ws.Range("N11").Value = 0.06666666666666667
$ws.Range("O11").Value = 1.0
$ws.Range("D12").Value = 0.16666666666666666
$ws.Range("E12").Value = 1.0
$ws.Range("F12").Value = 0.16666666666666666
$ws.Range("G12").Value = 1.0
$ws.Range("H12").Value = 0.16666666666666666
$ws.Range("I12").Value = 1.0
$ws.Range("F13").Value = 0.047619047619047616
$ws.Range("G13").Value = 1.0
$ws.Range("H13").Value = 0.09523809523809523
$ws.Range("I13").Value = 2.0
$ws.Range("J13").Value = 0.09523809523809523
$ws.Range("K13").Value = 2.0
$ws.Range("L13").Value = 0.09523809523809523
$ws.Range("M13").Value = 2.0
$ws.Range("N13").Value = 0.14285714285714285
$ws.Range("O13").Value = 3.0
$ws.Range("D14").Value = 0.043478260869565216
$ws.Range("E14").Value = 1.0
$ws.Range("F14").Value = 0.13043478260869565
$ws.Range("G14").Value = 3.0
$ws.Range("H14").Value = 0.17391304347826086
$ws.Range("I14").Value = 4.0
$ws.Range("N14").Value = 0.043478260869565216
$ws.Range("O14").Value = 1.0
$ws.Range("F15").Value = 0.03225806451612903
$ws.Range("G15").Value = 1.0
$ws.Range("H15").Value = 0.03225806451612903
$ws.Range("I15").Value = 1.0
$ws.Range("N15").Value = 0.03225806451612903
$ws.Range("O15").Value = 1.0
$ws.Range("D16").Value = 0.06451612903225806
$ws.Range("E16").Value = 2.0
$ws.Range("F16").Value = 0.06451612903225806
$ws.Range("G16").Value = 2.0
$ws.Range("H16").Value = 0.12903225806451613
$ws.Range("I16").Value = 4.0
$ws.Range("J16").Value = 0.06451612903225806
$ws.Range("K16").Value = 2.0
$ws.Range("L16").Value = 0.12903225806451613
$ws.Range("M16").Value = 4.0
$ws.Range("D19").Value = 0.10526315789473684
$ws.Range("E19").Value = 2.0
$ws.Range("F19").Value = 0.21052631578947367
$ws.Range("G19").Value = 4.0
$ws.Range("J19").Value = 0.05263157894736842
$ws.Range("K19").Value = 1.0
$ws.Range("L19").Value = 0.05263157894736842
$ws.Range("M19").Value = 1.0
$ws.Range("N19").Value = 0.21052631578947367
$ws.Range("O19").Value = 4.0
$ws.Range("H20").Value = 0.07142857142857142
$ws.Range("I20").Value = 2.0
$ws.Range("J20").Value = 0.10714285714285714
$ws.Range("K20").Value = 3.0
$ws.Range("D22").Value = 0.05263157894736842
$ws.Range("E22").Value = 1.0
$ws.Range("F22").Value = 0.10526315789473684
$ws.Range("G22").Value = 2.0
$ws.Range("L24").Value = 0.125
$ws.Range("M24").Value = 1.0
$ws.Range("N24").Value = 0.25
$ws.Range("O24").Value = 2.0
$ws.Range("D29").Value = 0.06896551724137931
$ws.Range("E29").Value = 2.0
$ws.Range("F29").Value = 0.13793103448275862
$ws.Range("G29").Value = 4.0
$ws.Range("L29").Value = 0.034482758620689655
$ws.Range("M29").Value = 1.0
$ws.Range("N29").Value = 0.06896551724137931
$ws.Range("O29").Value = 2.0
$ws.Range("D32").Value = 0.017857142857142856
$ws.Range("E32").Value = 1.0
$ws.Range("F32").Value = 0.07142857142857142
$ws.Range("G32").Value = 4.0
$ws.Range("J32").Value = 0.017857142857142856
$ws.Range("K32").Value = 1.0
$ws.Range("L32").Value = 0.07142857142857142
$ws.Range("M32").Value = 4.0
$ws.Range("D35").Value = 0.2
$ws.Range("E35").Value = 1.0
$ws.Range("F35").Value = 0.2
$ws.Range("G35").Value = 1.0
$ws.Range("H35").Value = 0.4
$ws.Range("I35").Value = 2.0
$ws.Range("H37").Value = 0.05555555555555555
$ws.Range("I37").Value = 1.0
$ws.Range("H38").Value = 0.07142857142857142
$ws.Range("I38").Value = 1.0
$ws.Range("J38").Value = 0.07142857142857142
$ws.Range("K38").Value = 1.0
$ws.Range("L38").Value = 0.07142857142857142
$ws.Range("M38").Value = 1.0
$ws.Range("N38").Value = 0.14285714285714285
$ws.Range("O38").Value = 2.0
$ws.Range("F41").Value = 0.02857142857142857
$ws.Range("G41").Value = 1.0
$ws.Range("H41").Value = 0.11428571428571428
$ws.Range("I41").Value = 4.0
$ws.Range("J41").Value = 0.02857142857142857
$ws.Range("K41").Value = 1.0
$ws.Range("H43").Value = 0.03125
$ws.Range("I43").Value = 1.0
$ws.Range("J43").Value = 0.09375
$ws.Range("K43").Value = 3.0
$ws.Range("D44").Value = 0.08333333333333333
$ws.Range("E44").Value = 2.0
$ws.Range("F44").Value = 0.08333333333333333
$ws.Range("G44").Value = 2.0
$ws.Range("H44").Value = 0.08333333333333333
$ws.Range("I44").Value = 2.0
$ws.Range("J44").Value = 0.041666666666666664
$ws.Range("K44").Value = 1.0
$ws.Range("L44").Value = 0.125
$ws.Range("M44").Value = 3.0
$ws.Range("D45").Value = 0.12
$ws.Range("E45").Value = 3.0
$ws.Range("F45").Value = 0.16
$ws.Range("G45").Value = 4.0
$ws.Range("J45").Value = 0.08
$ws.Range("K45").Value = 2.0
$ws.Range("L45").Value = 0.12
$ws.Range("M45").Value = 3.0
$ws.Range("N45").Value = 0.16
$ws.Range("O45").Value = 4.0
$ws.Range("D46").Value = 0.047619047619047616
$ws.Range("E46").Value = 1.0
$ws.Range("F46").Value = 0.09523809523809523
$ws.Range("G46").Value = 2.0
$ws.Range("H46").Value = 0.19047619047619047
$ws.Range("I46").Value = 4.0
$ws.Range("L46").Value = 0.047619047619047616
$ws.Range("M46").Value = 1.0
$ws.Range("N46").Value = 0.047619047619047616
$ws.Range("O46").Value = 1.0
$ws.Range("N47").Value = 0.027777777777777776
$ws.Range("O47").Value = 1.0
$ws.Range("D48").Value = 0.020833333333333332
$ws.Range("E48").Value = 1.0
$ws.Range("F48").Value = 0.041666666666666664
$ws.Range("G48").Value = 2.0
$ws.Range("H48").Value = 0.08333333333333333
$ws.Range("I48").Value = 4.0
$ws.Range("J48").Value = 0.08333333333333333
$ws.Range("K48").Value = 4.0
$ws.Range("D49").Value = 0.037037037037037035
$ws.Range("E49").Value = 2.0
$ws.Range("F49").Value = 0.037037037037037035
$ws.Range("G49").Value = 2.0
$ws.Range("F50").Value = 0.024390243902439025
$ws.Range("G50").Value = 1.0
$ws.Range("H50").Value = 0.04878048780487805
$ws.Range("I50").Value = 2.0
$ws.Range("J50").Value = 0.04878048780487805
$ws.Range("K50").Value = 2.0
$ws.Range("D51").Value = 0.06451612903225806
$ws.Range("E51").Value = 2.0
$ws.Range("F51").Value = 0.0967741935483871
$ws.Range("G51").Value = 3.0
$ws.Range("J51").Value = 0.0967741935483871
$ws.Range("K51").Value = 3.0
$ws.Range("L51").Value = 0.0967741935483871
$ws.Range("M51").Value = 3.0
$ws.Range("N51").Value = 0.12903225806451613
$ws.Range("O51").Value = 4.0
$ws.Range("D52").Value = 0.058823529411764705
$ws.Range("E52").Value = 1.0
$ws.Range("F52").Value = 0.058823529411764705
$ws.Range("G52").Value = 1.0
$ws.Range("H52").Value = 0.058823529411764705
$ws.Range("I52").Value = 1.0
$ws.Range("J52").Value = 0.11764705882352941
$ws.Range("K52").Value = 2.0
$ws.Range("L52").Value = 0.17647058823529413
$ws.Range("M52").Value = 3.0
$ws.Range("N52").Value = 0.17647058823529413
$ws.Range("O52").Value = 3.0
$ws.Range("F54").Value = 0.06451612903225806
$ws.Range("G54").Value = 2.0
$ws.Range("H54").Value = 0.06451612903225806
$ws.Range("I54").Value = 2.0
$ws.Range("L54").Value = 0.06451612903225806
$ws.Range("M54").Value = 2.0
$ws.Range("N54").Value = 0.12903225806451613
$ws.Range("O54").Value = 4.0
